$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("failing testcases")

# --- Fill new rows 29-32, columns A-D first (matches shared-string insertion order) ---
$ws.Cells.Item(29,1).Value = "small poster for Blade Runner"
$ws.Cells.Item(29,2).Value = "shows normal poster"
$ws.Cells.Item(29,3).Value = "is there a small poster?"

$ws.Cells.Item(30,1).Value = "genre of Blade Runner"
$ws.Cells.Item(30,2).Value = "no response"
$ws.Cells.Item(30,3).Value = "missing training item"
$ws.Cells.Item(30,4).Value = "fixed"

$ws.Cells.Item(31,1).Value = "cast of Blade Runner"
$ws.Cells.Item(31,2).Value = "error message"
$ws.Cells.Item(31,3).Value = "missing training item"
$ws.Cells.Item(31,4).Value = "fixed"

$ws.Cells.Item(32,1).Value = "Star Wars movies"
$ws.Cells.Item(32,2).Value = "error message"
$ws.Cells.Item(32,3).Value = "missing training item?"
$ws.Cells.Item(32,4).Value = "fixed"

# --- New "priority" column header ---
$ws.Cells.Item(1,5).Value = "priority"

# --- Fill the new priority column for rows 29-32 ---
$ws.Cells.Item(29,5).Value = "December milestone"
$ws.Cells.Item(30,5).Value = "December milestone"
$ws.Cells.Item(31,5).Value = "December milestone"
$ws.Cells.Item(32,5).Value = "December milestone"

# --- New row 33 (A, B, C, E only - no D) ---
$ws.Cells.Item(33,1).Value = "year of Star Wars"
$ws.Cells.Item(33,2).Value = "list years of all"
$ws.Cells.Item(33,3).Value = "ambiguous - need to clarify intent of request"
$ws.Cells.Item(33,5).Value = "December milestone"

# --- Update selection to match the new active cell ---
$ws.Range("F30").Select()
